# Added more text events
# Added hiding/showing the text box.
# Added changing the talking sfx.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: new lone numeric marker row (was START_SCENE/CUTSCENE/theme, now moved to row 2)
$ws.Range("A1").Value = 1
$ws.Range("B1").Value = ""
$ws.Range("C1").Value = ""

# Row 2: START_SCENE / CUTSCENE / theme (shifted down from old row 1)
$ws.Range("A2").Value = "START_SCENE"
$ws.Range("B2").Value = "CUTSCENE"
$ws.Range("C2").Value = "jazzy_retro_battle_theme"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = ""

# Row 3: Tanuki dialogue about disappearing (hide/show text box + talk sfx reset)
$ws.Range("A3").Value = "DIALOGUE"
$ws.Range("B3").Value = "Tanuki"
$ws.Range("C3").Value = "Im going to disappear[hide-text-box=t][set-talk-sfx=_] And then will you see me now?[hide-text-box=f]"
$ws.Range("D3").Value = "jazzy_retro_battle_theme"
$ws.Range("E3").Value = "tanuki_mario"
$ws.Range("F3").Value = "0,1"
$ws.Range("G3").Value = "END_DIALOGUE"

# Row 4: Frog dialogue that sets the talking sfx to take_damage
$ws.Range("A4").Value = "DIALOGUE"
$ws.Range("B4").Value = "Frog"
$ws.Range("C4").Value = "[set-talk-sfx=take_damage]Hmm"
$ws.Range("D4").Value = "frogs"
$ws.Range("E4").Value = "frog_mario"
$ws.Range("F4").Value = "LEFT"
$ws.Range("G4").Value = "END_DIALOGUE"
$ws.Range("H4").Value = ""
$ws.Range("I4").Value = ""

# Row 5: Tanuki dialogue "Hmmmmm"
$ws.Range("A5").Value = "DIALOGUE"
$ws.Range("B5").Value = "Tanuki"
$ws.Range("C5").Value = "Hmmmmm"
$ws.Range("D5").Value = "jazzy_retro_battle_theme"
$ws.Range("E5").Value = "tanuki_mario"
$ws.Range("F5").Value = "RIGHT"
$ws.Range("G5").Value = "frog_mario"
$ws.Range("H5").Value = "1, 1"
$ws.Range("I5").Value = "END_DIALOGUE"

# Row 6: END_SCENE (previously the old Tanuki BOO dialogue row; now cleared of extras)
$ws.Range("A6").Value = "END_SCENE"
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = ""
$ws.Range("G6").Value = ""
$ws.Range("H6").Value = ""
$ws.Range("I6").Value = ""

# Row 7: new END_GAME row (previously END_SCENE)
$ws.Range("A7").Value = "END_GAME"
$ws.Range("B7").Value = ""
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = ""
$ws.Range("G7").Value = ""
$ws.Range("H7").Value = ""
$ws.Range("I7").Value = ""

# Move the active selection to C3, matching the saved view state.
$ws.Range("C3").Select() | Out-Null
